$d = $word.ActiveDocument

$old = "is after using our product, users will never need to remember what they gonna do in the future. "
$new = "is to allow user to keep track of their task and also to increase their workflow by having simplicity and intuitive interface."

$range = $d.Content
$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
